# Apply "repull data, push all data, mean calculation" update:
# updates the dSF (column F) values for a set of rows to match
# freshly re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 3
    6  = -1
    7  = 1
    12 = 5
    15 = -4
    21 = -6
    23 = 0
    24 = 0
    29 = -8
    31 = -4
    32 = -5
    37 = 0
    38 = 1
    39 = -8
    46 = -6
    47 = -4
    48 = -2
    52 = 2
    54 = -1
    58 = -9
    60 = 0
    64 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
